$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing last header cell (H1) onto the two
# new header cells so they pick up the same style (bold, border, centered)
# that is used by the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for row 2
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 5
